$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-12-09T12:42:49.728Z","test@test.com","https://pnu.edu.sa","No","٩‏/١٢‏/٢٠٢٥","٣:٤٢:٤٩ م"),
    @("2025-12-09T12:55:50.316Z","rousol@ramworld.net","https://care.med.sa/","No","٩‏/١٢‏/٢٠٢٥","٣:٥٥:٥٠ م"),
    @("2025-12-09T12:58:11.542Z","rousol@ramworld.net","https://care.med.sa/","No","٩‏/١٢‏/٢٠٢٥","٣:٥٨:١١ م"),
    @("2025-12-09T13:03:53.508Z","rousol@ramworld.net","https://care.med.sa/","No","٩‏/١٢‏/٢٠٢٥","٤:٠٣:٥٣ م"),
    @("2025-12-09T13:08:16.699Z","rousol@ramworld.net","https://care.med.sa/","No","٩‏/١٢‏/٢٠٢٥","٤:٠٨:١٦ م"),
    @("2025-12-09T13:10:51.186Z","rousol@ramworld.net","https://care.med.sa/","Yes","٩‏/١٢‏/٢٠٢٥","٤:١٠:٥١ م"),
    @("2025-12-09T13:12:28.696Z","rousol@ramworld.net","https://care.med.sa/","Yes","٩‏/١٢‏/٢٠٢٥","٤:١٢:٢٨ م"),
    @("2025-12-09T13:13:11.151Z","rousol@ramworld.net","https://dga.gov.sa","Yes","٩‏/١٢‏/٢٠٢٥","٤:١٣:١١ م"),
    @("2025-12-09T13:14:42.767Z","rousol@ramworld.net","https://dga.gov.sa","Yes","٩‏/١٢‏/٢٠٢٥","٤:١٤:٤٢ م"),
    @("2025-12-10T06:25:33.109Z","rousol@ramworld.net","https://dga.gov.sa","Yes","١٠‏/١٢‏/٢٠٢٥","٩:٢٥:٣٣ ص"),
    @("2025-12-10T06:43:43.388Z","rousol@ramworld.net","https://dga.gov.sa","Yes","١٠‏/١٢‏/٢٠٢٥","٩:٤٣:٤٣ ص"),
    @("2025-12-10T07:06:05.108Z","rousol@ramworld.net","https://pnu.edu.sa","Yes","١٠‏/١٢‏/٢٠٢٥","١٠:٠٦:٠٥ ص"),
    @("2025-12-10T07:11:56.072Z","rousol@ramworld.net","https://hrdf.org.sa","No","١٠‏/١٢‏/٢٠٢٥","١٠:١١:٥٦ ص"),
    @("2025-12-10T07:18:53.455Z","rousol@ramworld.net","https://hrdf.org.sa","No","١٠‏/١٢‏/٢٠٢٥","١٠:١٨:٥٣ ص"),
    @("2025-12-10T08:14:22.546Z","rousol@ramworld.net","https://hrdf.org.sa","No","١٠‏/١٢‏/٢٠٢٥","١١:١٤:٢٢ ص"),
    @("2025-12-10T08:21:12.958Z","rousol@ramworld.net","https://hrdf.org.sa","No","١٠‏/١٢‏/٢٠٢٥","١١:٢١:١٢ ص"),
    @("2025-12-10T08:42:01.182Z","rousol@ramworld.net","https://pnu.edu.sa","No","١٠‏/١٢‏/٢٠٢٥","١١:٤٢:٠١ ص"),
    @("2025-12-10T08:49:02.598Z","rousol@ramworld.net","https://dga.gov.sa","No","١٠‏/١٢‏/٢٠٢٥","١١:٤٩:٠٢ ص"),
    @("2025-12-10T08:55:39.678Z","rousol@ramworld.net","https://dga.gov.sa","No","١٠‏/١٢‏/٢٠٢٥","١١:٥٥:٣٩ ص"),
    @("2025-12-10T08:58:52.955Z","rousol@ramworld.net","https://dga.gov.sa","No","١٠‏/١٢‏/٢٠٢٥","١١:٥٨:٥٢ ص"),
)

$startRow = 47
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}
